# I Aint Building That - 3118060751 update (#340)
# - duplicate the existing translation sheet into a new "Main_231223" sheet
#   placed before the original sheet (renamed "231221")
# - add a new translation row (Taggerung_IAintBuildingThat_SearchLabel / Filter / 필터)
#   to the new sheet, matching the existing row layout/styling

$wb = $excel.ActiveWorkbook

# --- Rename the existing (only) sheet, then duplicate it in front of itself ---
$orig = $wb.Worksheets.Item(1)
$orig.Copy($orig)

# After Copy(Before:=orig), the fresh duplicate sits at index 1 and the
# original (now pushed back) sits at index 2. Re-fetch by position since
# worksheet variables here track slot position, not object identity.
$main = $wb.Worksheets.Item(1)
$old  = $wb.Worksheets.Item(2)

$main.Name = "Main_231223"
$old.Name  = "231221"

# --- Add the new "SearchLabel" translation row (row 4) on the new sheet ---
# Shared strings: 19 Keyed+Taggerung_IAintBuildingThat_SearchLabel
#                 20 Taggerung_IAintBuildingThat_SearchLabel
#                 21 Filter
#                 22 (Korean) Filter translation
$main.Range("A4").Value = "Keyed+Taggerung_IAintBuildingThat_SearchLabel"
$main.Range("B4").Value = "Keyed"
$main.Range("C4").Value = "Taggerung_IAintBuildingThat_SearchLabel"
$main.Range("D4").Value = "Filter"
$main.Range("E4").Value = "필터"

# Match the white-fill formatting used by the other rows' A:E columns
# (copy the look of row 2, which already carries that style).
$main.Range("A2:B2").Copy()
$main.Range("A4:B4").PasteSpecial(-4122)
$main.Range("E2").Copy()
$main.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Select E11 on the new front sheet, matching the authored file ---
$main.Activate()
$main.Range("E11").Select()
